$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = $ws.Range("A4").Value2
$ws.Range("C5").Value = $ws.Range("A5").Value2
